# Update res_bus/vm_pu.xlsx results for the "case with 380 kV" run.
# For each bus-result row (2-25), columns B:F and I:N are replaced with the
# newly computed per-unit voltage magnitudes; column G (slack bus, always 1)
# and column A (bus index) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.038450854315047; $bf[0,2]=1.041055307278392; $bf[0,3]=1.047484124626031; $bf[0,4]=1.059320034148477
$ws.Range("B2:F2").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042305354612622; $in[0,1]=1.043548223093067; $in[0,2]=1.043835906299398; $in[0,3]=1.050246633662113; $in[0,4]=1.062049860313014; $in[0,5]=1.045030181266513
$ws.Range("I2:N2").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.039307683498232; $bf[0,2]=1.041690810562869; $bf[0,3]=1.048280844263789; $bf[0,4]=1.060279266528236
$ws.Range("B3:F3").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042538592217509; $in[0,1]=1.044050315127925; $in[0,2]=1.044282419196137; $in[0,3]=1.050855245979404; $in[0,4]=1.062822931938304; $in[0,5]=1.045532986329652
$ws.Range("I3:N3").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.039862585097072; $bf[0,2]=1.042102349547758; $bf[0,3]=1.048797213501648; $bf[0,4]=1.060901062951296
$ws.Range("B4:F4").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042688430876998; $in[0,1]=1.04437501694122; $in[0,2]=1.04457097935509; $in[0,3]=1.051249240390243; $in[0,4]=1.063323644612872; $in[0,5]=1.045858149256766
$ws.Range("I4:N4").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.040095978204457; $bf[0,2]=1.042275436943798; $bf[0,3]=1.049014494198756; $bf[0,4]=1.061162729976496
$ws.Range("B5:F5").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042751163355636; $in[0,1]=1.044511476073568; $in[0,2]=1.044692201875289; $in[0,3]=1.051414917807071; $in[0,4]=1.063534258526472; $in[0,5]=1.045994802176733
$ws.Range("I5:N5").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.040135172488035; $bf[0,2]=1.042304503513013; $bf[0,3]=1.049050988209884; $bf[0,4]=1.061206680442756
$ws.Range("B6:F6").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042761681164768; $in[0,1]=1.044534385460274; $in[0,2]=1.044712550466389; $in[0,3]=1.051442738209381; $in[0,4]=1.063569628189262; $in[0,5]=1.046017744097396
$ws.Range("I6:N6").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.03986570326704; $bf[0,2]=1.042104662051809; $bf[0,3]=1.048800116036803; $bf[0,4]=1.060904558324397
$ws.Range("B7:F7").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.04268927013279; $in[0,1]=1.044376840495086; $in[0,2]=1.044572599484215; $in[0,3]=1.051251454015934; $in[0,4]=1.063326458399005; $in[0,5]=1.045859975400287
$ws.Range("I7:N7").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.038740324361691; $bf[0,2]=1.041270010151387; $bf[0,3]=1.047753205187524; $bf[0,4]=1.059643980605181
$ws.Range("B8:F8").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042384401873878; $in[0,1]=1.043717945126371; $in[0,2]=1.04398688201332; $in[0,3]=1.050452278675039; $in[0,4]=1.062311022457901; $in[0,5]=1.045200144324572
$ws.Range("I8:N8").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.036760970769249; $bf[0,2]=1.039801813511019; $bf[0,3]=1.045914907585024; $bf[0,4]=1.05743124862967
$ws.Range("B9:F9").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041838933530697; $in[0,1]=1.042555518758623; $in[0,2]=1.04295203895765; $in[0,3]=1.049045477304738; $in[0,4]=1.060525464097366; $in[0,5]=1.04403606717805
$ws.Range("I9:N9").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.035443979808601; $bf[0,2]=1.038824836899029; $bf[0,3]=1.044693838623091; $bf[0,4]=1.055961953241809
$ws.Range("B10:F10").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.04146978409276; $in[0,1]=1.041779715791726; $in[0,2]=1.042260372932316; $in[0,3]=1.048108660024862; $in[0,4]=1.059337714098422; $in[0,5]=1.043259162481954
$ws.Range("I10:N10").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034874337449394; $bf[0,2]=1.03840224729428; $bf[0,3]=1.044166181224962; $bf[0,4]=1.055327142667858
$ws.Range("B11:F11").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041308642062667; $in[0,1]=1.041443596713364; $in[0,2]=1.041960468105064; $in[0,3]=1.047703273062502; $in[0,4]=1.058824045608612; $in[0,5]=1.042922566075946
$ws.Range("I11:N11").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034662841983486; $bf[0,2]=1.038245347581457; $bf[0,3]=1.043970348701412; $bf[0,4]=1.055091558361691
$ws.Range("B12:F12").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041248592461987; $in[0,1]=1.041318719518995; $in[0,2]=1.041849009797143; $in[0,3]=1.047552734811524; $in[0,4]=1.058633343146478; $in[0,5]=1.042797511541639
$ws.Range("I12:N12").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034708204154118; $bf[0,2]=1.038278999959449; $bf[0,3]=1.044012348043255; $bf[0,4]=1.055142082347418
$ws.Range("B13:F13").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041261482085038; $in[0,1]=1.041345507331976; $in[0,2]=1.041872920693238; $in[0,3]=1.047585023923204; $in[0,4]=1.058674245064044; $in[0,5]=1.042824337396387
$ws.Range("I13:N13").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034856853220658; $bf[0,2]=1.038389276500549; $bf[0,3]=1.044149990314287; $bf[0,4]=1.055307664848873
$ws.Range("B14:F14").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041303682302958; $in[0,1]=1.041433274882055; $in[0,2]=1.041951256150184; $in[0,3]=1.04769082868973; $in[0,4]=1.058808280098648; $in[0,5]=1.042912229586453
$ws.Range("I14:N14").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034948453456423; $bf[0,2]=1.038457230693022; $bf[0,3]=1.044234817898954; $bf[0,4]=1.05540971391519
$ws.Range("B15:F15").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041329657528047; $in[0,1]=1.041487347742619; $in[0,2]=1.041999513253737; $in[0,3]=1.047756023905813; $in[0,4]=1.058890876389527; $in[0,5]=1.042966379236682
$ws.Range("I15:N15").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.035481798350725; $bf[0,2]=1.038852892350344; $bf[0,3]=1.044728880264959; $bf[0,4]=1.056004113252341
$ws.Range("B16:F16").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041480451274757; $in[0,1]=1.041802018974816; $in[0,2]=1.042280268109671; $in[0,3]=1.04813556983358; $in[0,4]=1.059371818090562; $in[0,5]=1.043281497338122
$ws.Range("I16:N16").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.035816519411768; $bf[0,2]=1.039101201463792; $bf[0,3]=1.045039081106495; $bf[0,4]=1.056377341326158
$ws.Range("B17:F17").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041574693232454; $in[0,1]=1.041999353515967; $in[0,2]=1.042456269541674; $in[0,3]=1.048373719780792; $in[0,4]=1.059673671184016; $in[0,5]=1.043479112116956
$ws.Range("I17:N17").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.036011816599787; $bf[0,2]=1.039246079058264; $bf[0,3]=1.045220119478088; $bf[0,4]=1.056595174303212
$ws.Range("B18:F18").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041629537696426; $in[0,1]=1.042114436883933; $in[0,2]=1.04255888867327; $in[0,3]=1.048512653721881; $in[0,4]=1.059849798165522; $in[0,5]=1.043594358916504
$ws.Range("I18:N18").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.03607841803774; $bf[0,2]=1.0392954858335; $bf[0,3]=1.04528186638446; $bf[0,4]=1.056669472627141
$ws.Range("B19:F19").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041648216970645; $in[0,1]=1.04215367416068; $in[0,2]=1.042593872406653; $in[0,3]=1.048560030824633; $in[0,4]=1.059909863247025; $in[0,5]=1.043633651914686
$ws.Range("I19:N19").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.03578060078762; $bf[0,2]=1.039074555777243; $bf[0,3]=1.045005788777764; $bf[0,4]=1.056337283473473
$ws.Range("B20:F20").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041564594914945; $in[0,1]=1.04197818329584; $in[0,2]=1.042437390333898; $in[0,3]=1.048348165939135; $in[0,4]=1.059641278865144; $in[0,5]=1.043457911832688
$ws.Range("I20:N20").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.03481307713983; $bf[0,2]=1.038356800893298; $bf[0,3]=1.044109453585299; $bf[0,4]=1.055258899046866
$ws.Range("B21:F21").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.04129126074779; $in[0,1]=1.041407430278974; $in[0,2]=1.041928189966148; $in[0,3]=1.047659670694288; $in[0,4]=1.05876880743535; $in[0,5]=1.042886348281072
$ws.Range("I21:N21").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034205307364727; $bf[0,2]=1.037905919240487; $bf[0,3]=1.043546835539509; $bf[0,4]=1.054582106856304
$ws.Range("B22:F22").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041118281044436; $in[0,1]=1.041048416063657; $in[0,2]=1.041607687155015; $in[0,3]=1.047227021463325; $in[0,4]=1.05822081130108; $in[0,5]=1.042526824224389
$ws.Range("I22:N22").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.034527444946681; $bf[0,2]=1.038144901744029; $bf[0,3]=1.043845000032211; $bf[0,4]=1.054940769969516
$ws.Range("B23:F23").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041210087128492; $in[0,1]=1.041238750997883; $in[0,2]=1.041777624400652; $in[0,3]=1.047456354295444; $in[0,4]=1.058511260725942; $in[0,5]=1.042717429456054
$ws.Range("I23:N23").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.03579683068236; $bf[0,2]=1.039086595683102; $bf[0,3]=1.045020831828561; $bf[0,4]=1.056355383474791
$ws.Range("B24:F24").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.041569158297108; $in[0,1]=1.04198774926444; $in[0,2]=1.042445921157504; $in[0,3]=1.048359712542802; $in[0,4]=1.059655915370863; $in[0,5]=1.043467491386061
$ws.Range("I24:N24").Value = $in

$bf = New-Object "object[,]" 1,5
$bf[0,0]=1.02; $bf[0,1]=1.037272232956973; $bf[0,2]=1.040181063527504; $bf[0,3]=1.046389372100772; $bf[0,4]=1.058002267426352
$ws.Range("B25:F25").Value = $bf

$in = New-Object "object[,]" 1,6
$in[0,0]=1.042856189004508; $in[0,1]=1.042856189004508; $in[0,2]=1.043219887728785; $in[0,3]=1.049408989693935; $in[0,4]=1.060986617982911; $in[0,5]=1.044337164410171
$ws.Range("I25:N25").Value = $in

Write-Output "case with 380 kV done"
